$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C, shifting old C..I to D..J
$ws.Columns("C").Insert()

# Header for the new Industry column
$ws.Range("C1").Value = "Industry"

# Industry values per row (2-65)
$industries = @{
    2 = "Finance"
    3 = "Retailing"
    4 = "Banks"
    5 = "Banks"
    6 = "Capital Markets"
    7 = "Healthcare Services"
    8 = "Capital Markets"
    9 = "Capital Markets"
    10 = "Automobiles"
    11 = "Banks"
    12 = "Pharmaceuticals & Biotechnology"
    13 = "IT - Software"
    14 = "Healthcare Services"
    15 = "Transport Services"
    16 = "Realty"
    17 = "Healthcare Services"
    18 = "Cement & Cement Products"
    19 = "Transport Services"
    20 = "Retailing"
    21 = "Insurance"
    22 = "Retailing"
    23 = "Financial Technology (Fintech)"
    24 = "IT - Software"
    25 = "Automobiles"
    26 = "IT - Software"
    27 = "Construction"
    28 = "Cement & Cement Products"
    29 = "Healthcare Services"
    30 = "Realty"
    31 = "Consumer Durables"
    32 = "Consumer Durables"
    33 = "Aerospace & Defense"
    34 = "Retailing"
    35 = "Consumer Durables"
    36 = "Auto Components"
    37 = "Finance"
    38 = "Finance"
    39 = "Banks"
    40 = "Industrial Products"
    41 = "Leisure Services"
    42 = "Industrial Products"
    43 = "IT - Software"
    44 = "Electrical Equipment"
    45 = "Banks"
    46 = "Finance"
    47 = "Pharmaceuticals & Biotechnology"
    48 = "Pharmaceuticals & Biotechnology"
    49 = "Consumer Durables"
    50 = "Finance"
    51 = "Consumer Durables"
    52 = "Banks"
    53 = "Paper, Forest & Jute Products"
    54 = "Industrial Manufacturing"
    55 = "Retailing"
    56 = "Industrial Manufacturing"
    57 = "Finance"
    58 = "Personal Products"
    59 = "Commercial Services & Supplies"
    60 = "Automobiles"
    61 = "Leisure Services"
    62 = "Industrial Manufacturing"
    63 = "Capital Markets"
    64 = "Finance"
    65 = "Cement & Cement Products"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item([int]$row, 3).Value = $industries[$row]
}

